# Atualiza a planilha de previsao do tempo de Sao Paulo com os dados mais
# recentes coletados pelo bot: a janela de previsao "andou" um dia para a
# frente (a linha que era "ter. 24" sai da tabela) e uma nova linha (dia
# "sex. 04") e adicionada ao final, alem da atualizacao dos valores de
# temperatura/umidade/indice UV de cada dia que permanece na tabela.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dia | Temp. maxima | Temp. minima | Umidade do dia | Umidade da noite | Indice UV do dia | Indice UV da noite
$dados = @{
    2  = @("qua. 25", "34°", "21°", "39%", "55%", "10 de 11", "0 de 11")
    3  = @("qui. 26", "36°", "22°", "31%", "48%", "10 de 11", "0 de 11")
    4  = @("sex. 27", "29°", "16°", "60%", "84%", "10 de 11", "0 de 11")
    5  = @("sáb. 28", "21°", "15°", "72%", "85%", "7 de 11",  "0 de 11")
    6  = @("dom. 29", "26°", "16°", "67%", "86%", "10 de 11", "0 de 11")
    7  = @("seg. 30", "32°", "19°", "53%", "67%", "10 de 11", "0 de 11")
    8  = @("ter. 01", "35°", "21°", "36%", "52%", "Extremo",  "0 de 11")
    9  = @("qua. 02", "36°", "20°", "34%", "61%", "Extremo",  "0 de 11")
    10 = @("qui. 03", "32°", "19°", "48%", "76%", "Extremo",  "0 de 11")
    11 = @("sex. 04", "29°", "19°", "58%", "79%", "Extremo",  "0 de 11")
}

$colunas = @("A", "B", "C", "D", "E", "F", "G")

foreach ($linha in ($dados.Keys | Sort-Object)) {
    $valores = $dados[$linha]
    for ($i = 0; $i -lt $colunas.Length; $i++) {
        $celula = $ws.Range($colunas[$i] + $linha)
        $valor = $valores[$i]

        if ($valor.EndsWith("%")) {
            # Os percentuais de umidade sao textos (ex.: "39%"), nao numeros -
            # usamos o prefixo de aspa simples para evitar que o Excel
            # converta o texto digitado em um numero percentual, e depois
            # limpamos a formatacao de "texto forcado" que isso aplica para
            # manter a celula com o estilo padrao da tabela.
            $celula.Value = "'" + $valor
            $celula.ClearFormats()
        } else {
            $celula.Value = $valor
        }
    }
}
